$d = $word.ActiveDocument

# Locate "Coler" (the typo inside "Semester Project: Swamp Coler") so we
# don't rely on hard-coded character offsets.
$findRng = $d.Content
$found = $findRng.Find.Execute("Coler", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'Coler'"
}

# $findRng now spans "Coler"; the split point is right after the leading "C",
# i.e. right before "oler".
$splitPos = $findRng.Start + 1

# Insert the missing "o" so the word becomes "Cooler".
$insertRng = $d.Range($splitPos, $splitPos)
$insertRng.InsertBefore("o")

# Toggling a formatting property on just the newly inserted character keeps
# it from being silently re-merged with its neighboring runs, so the
# paragraph ends up as three runs: "Semester Project: Swamp C" / "o" / "oler"
# (matching formatting, just split across separate <w:r> elements).
$newCharRng = $d.Range($splitPos, $splitPos + 1)
$newCharRng.Bold = 1
$newCharRng.Bold = 0
